$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Photometric-Opt")
$ws.Range("D4").Value = 9.919356917703732
$ws.Range("K4").Value = 4.304025912258477
$ws.Range("D5").Value = 9.886936539789014
$ws.Range("K5").Value = 3.981254205101203
$ws.Range("B6").Value = 9.919356917703732
$ws.Range("C6").Value = 9.886936539789014
$ws.Range("E6").Value = 7.117976427563766
$ws.Range("F6").Value = 9.99217094197501
$ws.Range("I6").Value = 4.304025912258477
$ws.Range("J6").Value = 3.981254205101203
$ws.Range("L6").Value = 5.106398705457906
$ws.Range("M6").Value = 3.930611463694983
$ws.Range("D7").Value = 7.117976427563766
$ws.Range("K7").Value = 5.106398705457906
$ws.Range("D8").Value = 9.99217094197501
$ws.Range("K8").Value = 3.930611463694983
$ws.Range("D14").Value = 9.919356917702338
$ws.Range("K14").Value = 4.264097540967208
$ws.Range("D15").Value = 9.886936539793002
$ws.Range("K15").Value = 4.265870609220072
$ws.Range("B16").Value = 9.919356917702338
$ws.Range("C16").Value = 9.886936539793002
$ws.Range("E16").Value = 7.117976427563251
$ws.Range("F16").Value = 9.992170941966545
$ws.Range("I16").Value = 4.264097540967208
$ws.Range("J16").Value = 4.265870609220072
$ws.Range("L16").Value = 4.269155611861901
$ws.Range("M16").Value = 4.26838847448473
$ws.Range("D17").Value = 7.117976427563251
$ws.Range("K17").Value = 4.269155611861901
$ws.Range("D18").Value = 9.992170941966545
$ws.Range("K18").Value = 4.26838847448473

$ws = $wb.Worksheets.Item("Photometric-Pess")
$ws.Range("D4").Value = 10.84289138992081
$ws.Range("K4").Value = 2.614010165875955
$ws.Range("D5").Value = 10.8728026011959
$ws.Range("K5").Value = 2.625541679111308
$ws.Range("B6").Value = 10.84289138992081
$ws.Range("C6").Value = 10.8728026011959
$ws.Range("E6").Value = 10.5527857504133
$ws.Range("F6").Value = 10.84918233753177
$ws.Range("I6").Value = 2.614010165875955
$ws.Range("J6").Value = 2.625541679111308
$ws.Range("L6").Value = 4.026938695190727
$ws.Range("M6").Value = 2.617148387586216
$ws.Range("D7").Value = 10.5527857504133
$ws.Range("K7").Value = 4.026938695190727
$ws.Range("D8").Value = 10.84918233753177
$ws.Range("K8").Value = 2.617148387586216
$ws.Range("D14").Value = 10.84289138992082
$ws.Range("K14").Value = 3.453247927247328
$ws.Range("D15").Value = 10.87280260119602
$ws.Range("K15").Value = 3.453996048594282
$ws.Range("B16").Value = 10.84289138992082
$ws.Range("C16").Value = 10.87280260119602
$ws.Range("E16").Value = 10.55278575041431
$ws.Range("F16").Value = 10.84918233753237
$ws.Range("I16").Value = 3.453247927247328
$ws.Range("J16").Value = 3.453996048594282
$ws.Range("L16").Value = 3.453945669386608
$ws.Range("M16").Value = 3.454313084573598
$ws.Range("D17").Value = 10.55278575041431
$ws.Range("K17").Value = 3.453945669386608
$ws.Range("D18").Value = 10.84918233753237
$ws.Range("K18").Value = 3.454313084573598

$ws = $wb.Worksheets.Item("Spectroscopic-Opt")
$ws.Range("D4").Value = 9.041864007425653
$ws.Range("K4").Value = 12.42479335229207
$ws.Range("D5").Value = 9.073006014803994
$ws.Range("K5").Value = 12.42567251460646
$ws.Range("B6").Value = 9.041864007425653
$ws.Range("C6").Value = 9.073006014803994
$ws.Range("E6").Value = 9.703072181449262
$ws.Range("F6").Value = 8.960981046584797
$ws.Range("I6").Value = 12.42479335229207
$ws.Range("J6").Value = 12.42567251460646
$ws.Range("L6").Value = 13.13783480376886
$ws.Range("M6").Value = 12.38050782994369
$ws.Range("D7").Value = 9.703072181449262
$ws.Range("K7").Value = 13.13783480376886
$ws.Range("D8").Value = 8.960981046584797
$ws.Range("K8").Value = 12.38050782994369
$ws.Range("D14").Value = 12.08381760429481
$ws.Range("K14").Value = 0.474316134857282
$ws.Range("D15").Value = 12.04577796566226
$ws.Range("K15").Value = 0.465363295794745
$ws.Range("B16").Value = 12.08381760429481
$ws.Range("C16").Value = 12.04577796566226
$ws.Range("E16").Value = 12.796691355624
$ws.Range("F16").Value = 12.09867764770548
$ws.Range("I16").Value = 0.474316134857282
$ws.Range("J16").Value = 0.465363295794745
$ws.Range("L16").Value = 0.4687308738595936
$ws.Range("M16").Value = 0.4701064277266727
$ws.Range("D17").Value = 12.796691355624
$ws.Range("K17").Value = 0.4687308738595936
$ws.Range("D18").Value = 12.09867764770548
$ws.Range("K18").Value = 0.4701064277266727

$ws = $wb.Worksheets.Item("Spectroscopic-Pess")
$ws.Range("D4").Value = 8.062588612130609
$ws.Range("K4").Value = 10.48676195248214
$ws.Range("D5").Value = 8.006569333036184
$ws.Range("K5").Value = 10.47856220274357
$ws.Range("B6").Value = 8.062588612130609
$ws.Range("C6").Value = 8.006569333036184
$ws.Range("E6").Value = 8.856491601352039
$ws.Range("F6").Value = 7.859104476890197
$ws.Range("I6").Value = 10.48676195248214
$ws.Range("J6").Value = 10.47856220274357
$ws.Range("L6").Value = 11.20511984687745
$ws.Range("M6").Value = 10.43795098231388
$ws.Range("D7").Value = 8.856491601352039
$ws.Range("K7").Value = 11.20511984687745
$ws.Range("D8").Value = 7.859104476890197
$ws.Range("K8").Value = 10.43795098231388
$ws.Range("D14").Value = 11.69423127716002
$ws.Range("K14").Value = 0.3068155992074563
$ws.Range("D15").Value = 11.6573241284687
$ws.Range("K15").Value = 0.2983315187058641
$ws.Range("B16").Value = 11.69423127716002
$ws.Range("C16").Value = 11.6573241284687
$ws.Range("E16").Value = 12.47793868942139
$ws.Range("F16").Value = 11.72636687679834
$ws.Range("I16").Value = 0.3068155992074563
$ws.Range("J16").Value = 0.2983315187058641
$ws.Range("L16").Value = 0.3011553713424278
$ws.Range("M16").Value = 0.3025150120135459
$ws.Range("D17").Value = 12.47793868942139
$ws.Range("K17").Value = 0.3011553713424278
$ws.Range("D18").Value = 11.72636687679834
$ws.Range("K18").Value = 0.3025150120135459
